# Documentacion en mostrar inventario
# - Row 12 is replaced with a new inventory item ("Velador mediano de roble").
# - The old rows 13-22 (test/demo entries) are removed entirely, shrinking
#   the used range from A1:F22 down to A1:F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the data in row 12 with the new item.
$ws.Range("A12").Value = 24
$ws.Range("B12").Value = "Velador mediano de roble"
$ws.Range("C12").Value = 15
$ws.Range("D12").Value = 140
$ws.Range("E12").Value = "unidad"
$ws.Range("F12").Value = "Dormitorio"

# Remove the now-obsolete rows 13 through 22.
$ws.Range("A13:F22").EntireRow.Delete()
